# MirrorMe - Example.xlsx edit
#
# The single "TTexts" table (rows 1-9: 2 header rows + 7 data rows,
# TText_01..TText_07) is split into two separate tables:
#   - a "Template" table (new rows 1-7)   : TText_03 .. TText_07
#   - a "Binding" table  (new rows 10-13) : TText_01 .. TText_02
#
# Strategy: Range.Copy(Destination) so cell styles / number formats travel
# with the content (same as selecting rows and copy/pasting them
# interactively). The tricky bit is that two of the original 9 rows (8
# and 9) must disappear entirely (no leftover <row> element at all), and
# EntireRow.Delete() shifts everything below the deleted row up - so the
# Binding table (which ends up below where rows 8/9 used to be) is first
# written two rows "too low" (rows 12-15) and then lands on its real
# target (rows 10-13) once rows 8 and 9 are deleted at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Stage the content that will become the "Binding" table -----------
# Written 2 rows below its final destination; the later deletion of rows
# 8-9 shifts it up by 2 into its real place (rows 10-13).
$ws.Range("A1:F2").Copy($ws.Range("A12"))   # header -> future rows 10-11
$ws.Range("A3:G4").Copy($ws.Range("A14"))   # TText_01/02 -> future rows 12-13

# --- 2) Move the "Template" data rows up into rows 3-7 --------------------
# Each source row is read before anything overwrites it.
$ws.Range("A5:G5").Copy($ws.Range("A3"))    # TText_03 -> row 3 (was row 5)
$ws.Range("A6:G6").Copy($ws.Range("A4"))    # TText_04 -> row 4 (was row 6)
$ws.Range("A7:G7").Copy($ws.Range("A5"))    # TText_05 -> row 5 (was row 7)
$ws.Range("A8:G8").Copy($ws.Range("A6"))    # TText_06 -> row 6 (was row 8)
$ws.Range("A9:G9").Copy($ws.Range("A7"))    # TText_07 -> row 7 (was row 9)

# --- 3) Drop the now-redundant original rows 8 and 9 entirely ------------
# Delete from the bottom up so the row indices stay valid; this also
# shifts the staged rows 12-15 up to their real homes, rows 10-13.
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(8).Delete()

# --- 4) Clear stray cells left behind by the row copies above ------------
# Range.Copy only writes cells that existed in the source range, so cells
# that were non-empty at the destination before the copy (but empty in
# the source) are not cleared automatically.
$ws.Range("C3").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("G3").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("G4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("G6").ClearContents()
$ws.Range("C7").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("F7").ClearContents()

# --- 5) Fix up the text that is genuinely new rather than just moved -----

# Rows 1-2 stay in place and become the "Template" table header.
$ws.Range("A1").Value = "[Template]"
$ws.Range("A2").Value = "Template"
$ws.Range("C2").Value = "TTextID"

# Rows 10-11 (now in their final place) become the "Binding" table header.
$ws.Range("A10").Value = "[Binding]"
$ws.Range("A11").Value = "Binding"
$ws.Range("C11").Value = "TTextID"

# Row 5 (the TText_05 phrase, shifted up from row 7) gets an extended
# phrase.
$ws.Range("E5").Value = "De [bestuurder] is niet verkiesbaar tot lid van de ondernemingsraad van de [onderneming]."

# --- 6) Row heights --------------------------------------------------------
# Range.Copy does not bring row heights along, so set them explicitly to
# match the moved content.
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(7).RowHeight = 30
$ws.Rows.Item(12).RowHeight = 60
$ws.Rows.Item(13).RowHeight = 30

# --- 7) Match the saved selection / active cell in the edited workbook ---
$ws.Range("C8").Select()
